$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update F646 value
$ws.Range("F646").Value = 5566.35789109

# Row 647
$ws.Range("A646").Copy($ws.Range("A647"))
$ws.Range("A647").Value = 45116.41666666666
$ws.Range("B647").Value = 30299.25
$ws.Range("C647").Value = 30453.27
$ws.Range("D647").Value = 30080.24
$ws.Range("E647").Value = 30174.62
$ws.Range("F647").Value = 5874.45936717

# Row 648
$ws.Range("A647").Copy($ws.Range("A648"))
$ws.Range("A648").Value = 45117.41666666666
$ws.Range("B648").Value = 30175.34
$ws.Range("C648").Value = 31042.51
$ws.Range("D648").Value = 29965.03
$ws.Range("E648").Value = 30423.95
$ws.Range("F648").Value = 18369.45647798

# Row 649
$ws.Range("A648").Copy($ws.Range("A649"))
$ws.Range("A649").Value = 45118.41666666666
$ws.Range("B649").Value = 30422.95
$ws.Range("C649").Value = 30809.56
$ws.Range("D649").Value = 30320.36
$ws.Range("E649").Value = 30631.36
$ws.Range("F649").Value = 14390.16504579

# Row 650
$ws.Range("A649").Copy($ws.Range("A650"))
$ws.Range("A650").Value = 45119.41666666666
$ws.Range("B650").Value = 30633.89
$ws.Range("C650").Value = 30982
$ws.Range("D650").Value = 30227.25
$ws.Range("E650").Value = 30396.78
$ws.Range("F650").Value = 20184.77143358

# Row 651
$ws.Range("A650").Copy($ws.Range("A651"))
$ws.Range("A651").Value = 45120.41666666666
$ws.Range("B651").Value = 30395.64
$ws.Range("C651").Value = 31829
$ws.Range("D651").Value = 30258.46
$ws.Range("E651").Value = 31482.21
$ws.Range("F651").Value = 36831.45497786

# Row 652
$ws.Range("A651").Copy($ws.Range("A652"))
$ws.Range("A652").Value = 45121.41666666666
$ws.Range("B652").Value = 31483.23
$ws.Range("C652").Value = 31644.47
$ws.Range("D652").Value = 29940.08
$ws.Range("E652").Value = 30333.65
$ws.Range("F652").Value = 28504.11311169

# Row 653
$ws.Range("A652").Copy($ws.Range("A653"))
$ws.Range("A653").Value = 45122.41666666666
$ws.Range("B653").Value = 30332.66
$ws.Range("C653").Value = 30403.97
$ws.Range("D653").Value = 30267.04
$ws.Range("E653").Value = 30299
$ws.Range("F653").Value = 4039.37147264

# Row 654
$ws.Range("A653").Copy($ws.Range("A654"))
$ws.Range("A654").Value = 45123.41666666666
$ws.Range("B654").Value = 30300.6
$ws.Range("C654").Value = 30457.63
$ws.Range("D654").Value = 30078.23
$ws.Range("E654").Value = 30250.49
$ws.Range("F654").Value = 6357.2037676

# Row 655
$ws.Range("A654").Copy($ws.Range("A655"))
$ws.Range("A655").Value = 45124.41666666666
$ws.Range("B655").Value = 30248.97
$ws.Range("C655").Value = 30342.59
$ws.Range("D655").Value = 29678.15
$ws.Range("E655").Value = 30154.32
$ws.Range("F655").Value = 16010.77083874

# Row 656
$ws.Range("A655").Copy($ws.Range("A656"))
$ws.Range("A656").Value = 45125.41666666666
$ws.Range("B656").Value = 30152.07
$ws.Range("C656").Value = 30243.7
$ws.Range("D656").Value = 29522.25
$ws.Range("E656").Value = 29868.81
$ws.Range("F656").Value = 16104.96081001

# Row 657
$ws.Range("A656").Copy($ws.Range("A657"))
$ws.Range("A657").Value = 45126.41666666666
$ws.Range("B657").Value = 29863.81
$ws.Range("C657").Value = 30201.29
$ws.Range("D657").Value = 29770.34
$ws.Range("E657").Value = 29921.83
$ws.Range("F657").Value = 12551.08599458
